$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Plain text/volume updates: row -> @{ D = new price text (or $null); E = new volume text }
# D values that parse as plain numbers are handled separately further below so Excel
# doesn't silently coerce them from text into numeric cells.
$plainUpdates = @{
    2 = @{ D = '34.069.77'; E = '  -0.89%  ' }
    3 = @{ D = '1.788.44'; E = '  -1.56%  ' }
    4 = @{ D = $null; E = '  -0.05%  ' }
    5 = @{ D = $null; E = '  -1.65%  ' }
    6 = @{ D = $null; E = '  +1.79%  ' }
    7 = @{ D = $null; E = '  -0.05%  ' }
    8 = @{ D = $null; E = '  -0.40%  ' }
    9 = @{ D = $null; E = '  +1.23%  ' }
    10 = @{ D = $null; E = '  -0.60%  ' }
    11 = @{ D = $null; E = '  -2.81%  ' }
    12 = @{ D = $null; E = '  -0.63%  ' }
    13 = @{ D = '2.047.01'; E = '  -1.55%  ' }
    14 = @{ D = $null; E = '  +11.71%  ' }
    15 = @{ D = '1.785.44'; E = '  -1.05%  ' }
    16 = @{ D = $null; E = '  -1.28%  ' }
    17 = @{ D = '34.059.93'; E = '  -0.87%  ' }
    18 = @{ D = $null; E = '  -2.63%  ' }
    19 = @{ D = $null; E = '  -1.11%  ' }
    20 = @{ D = $null; E = '  -2.70%  ' }
    21 = @{ D = '0.0₃0742'; E = '  -0.98%  ' }
    22 = @{ D = $null; E = '  -0.01%  ' }
    23 = @{ D = $null; E = '  -0.54%  ' }
    24 = @{ D = $null; E = '  -2.05%  ' }
    25 = @{ D = $null; E = '  -2.49%  ' }
    26 = @{ D = $null; E = '  -2.89%  ' }
    27 = @{ D = $null; E = '  -1.06%  ' }
    28 = @{ D = $null; E = '  -1.46%  ' }
    29 = @{ D = $null; E = '  -2.01%  ' }
    31 = @{ D = $null; E = '  -0.04%  ' }
    32 = @{ D = $null; E = '  +0.17%  ' }
    33 = @{ D = $null; E = '  -1.15%  ' }
    34 = @{ D = $null; E = '  +1.27%  ' }
    35 = @{ D = $null; E = '  +1.11%  ' }
    36 = @{ D = '1.450.51'; E = '  -8.08%  ' }
    37 = @{ D = $null; E = '  +0.13%  ' }
    38 = @{ D = $null; E = '  +0.33%  ' }
    39 = @{ D = $null; E = '  -1.27%  ' }
    40 = @{ D = $null; E = '  -1.87%  ' }
    41 = @{ D = $null; E = '  -2.02%  ' }
    42 = @{ D = $null; E = '  -0.02%  ' }
    43 = @{ D = $null; E = '  -1.49%  ' }
    44 = @{ D = $null; E = '  -1.66%  ' }
    45 = @{ D = $null; E = '  -1.81%  ' }
    46 = @{ D = $null; E = '  +0.46%  ' }
    47 = @{ D = '1.947.21'; E = '  -1.47%  ' }
    48 = @{ D = $null; E = '  -0.06%  ' }
    49 = @{ D = $null; E = '  +0.03%  ' }
    50 = @{ D = $null; E = '  +7.74%  ' }
    51 = @{ D = $null; E = '  -3.36%  ' }
}

foreach ($row in $plainUpdates.Keys) {
    $entry = $plainUpdates[$row]
    if ($null -ne $entry.D) {
        $ws.Cells.Item($row, 4).Value = $entry.D
    }
    if ($null -ne $entry.E) {
        $ws.Cells.Item($row, 5).Value = $entry.E
    }
}

# D-column values that look like plain numbers (e.g. '226.78') must be forced to stay
# text (matching the original inlineStr cells) instead of being coerced to numeric.
# Route them through a text formula + paste-values so the stored cell type is string,
# without touching the cell's number format/style.
$numericLookingD = @{
    5 = '226.78'
    8 = '31.21'
    9 = '46.04'
    11 = '0.0660'
    14 = '11.39'
    16 = '0.635'
    18 = '4.22'
    19 = '69.47'
    20 = '253.28'
    24 = '4.29'
    26 = '156.58'
    32 = '0.0516'
    33 = '1.21'
    40 = '83.32'
    43 = '0.901'
    44 = '2.10'
    48 = '5.74'
    50 = '11.87'
    51 = '51.32'
}

foreach ($row in $numericLookingD.Keys) {
    $text = $numericLookingD[$row]
    $cell = $ws.Cells.Item($row, 4)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
